$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Entrada")

# Row 2 <-> Row 3 swap (labels swapped: DEVOLUÃÃO now in row 2, FERRAMENTAS/ MATRIZARIA now in row 3)
$ws.Range("A2").Value = "DEVOLUÃÃO"
$ws.Range("B2").Value = "R$ 653.089,99"
$ws.Range("D2").Value = "R$ 653.089,99"
$ws.Range("E2").Value = "R$ 653.089,99"
$ws.Range("F2").Value = "100,00 %"

$ws.Range("A3").Value = "FERRAMENTAS/ MATRIZARIA"
$ws.Range("B3").Value = "R$ 471.474,40"
$ws.Range("D3").Value = "R$ 471.474,40"
$ws.Range("E3").Value = "R$ 680.000,00"
$ws.Range("F3").Value = "69,33 %"

# Row 5: REFUGO REAL (PROCESSO) - values updated
$ws.Range("B5").Value = "R$ 309.142,60"
$ws.Range("D5").Value = "R$ 309.142,60"
$ws.Range("E5").Value = "R$ 309.142,60"

# Row 6 <-> Row 7 <-> Row 8 <-> Row 9 cyclic label shift
# Row 6 becomes CUSTO DESENVOLVIMENTO
$ws.Range("A6").Value = "CUSTO DESENVOLVIMENTO"
$ws.Range("B6").Value = "R$ 221.419,59"
$ws.Range("C6").Value = "R$ 0,00"
$ws.Range("D6").Value = "R$ 221.419,59"
$ws.Range("E6").Value = "R$ 221.419,59"
$ws.Range("F6").Value = "100,00 %"

# Row 7 becomes MANUTENCAO
$ws.Range("A7").Value = "MANUTENCAO"
$ws.Range("B7").Value = "R$ 205.025,51"
$ws.Range("C7").Value = "R$ 191.083,29"
$ws.Range("D7").Value = "R$ 396.108,80"
$ws.Range("E7").Value = "R$ 480.000,00"
$ws.Range("F7").Value = "82,52 %"

# Row 8 becomes FRETES
$ws.Range("A8").Value = "FRETES"
$ws.Range("B8").Value = "R$ 185.370,49"
$ws.Range("D8").Value = "R$ 185.370,49"
$ws.Range("E8").Value = "R$ 376.000,00"
$ws.Range("F8").Value = "49,30 %"

# Row 9 becomes REFUGO MP+CP*
$ws.Range("A9").Value = "REFUGO MP+CP*"
$ws.Range("B9").Value = "R$ 184.283,97"
$ws.Range("D9").Value = "R$ 184.283,97"
$ws.Range("E9").Value = "R$ 280.000,00"
$ws.Range("F9").Value = "65,82 %"

# Row 11: DESP. INDUSTRIAL - values updated
$ws.Range("B11").Value = "R$ 60.025,93"
$ws.Range("C11").Value = "R$ 104.730,86"
$ws.Range("D11").Value = "R$ 164.756,79"
$ws.Range("F11").Value = "36,61 %"

# Row 12: EMBALAGENS - values updated
$ws.Range("B12").Value = "R$ 58.564,60"
$ws.Range("D12").Value = "R$ 130.829,61"
$ws.Range("F12").Value = "81,77 %"

# Row 14: FERRAMENTARIA/MAN FR - values updated
$ws.Range("B14").Value = "R$ 11.114,15"
$ws.Range("D14").Value = "R$ 11.114,15"
$ws.Range("F14").Value = "31,75 %"

# Row 18: Total Geral - values updated
$ws.Range("B18").Value = "R$ 2.677.674,49"
$ws.Range("C18").Value = "R$ 481.379,59"
$ws.Range("D18").Value = "R$ 3.159.054,08"
$ws.Range("E18").Value = "R$ 4.261.954,26"
$ws.Range("F18").Value = "74,12 %"
